# The CI "Generate Report for Handoff" step flips the freshly-handed-off
# locale rows from "In Translation" to "Ready for handoff" and stamps the
# new handoff timestamps on the Overview summary sheet plus each per-locale
# detail sheet (zh-cn / de-de).
$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status + HO date (col G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-29 10:59:46"
# Status text grew longer ("In Translation" -> "Ready for handoff"), so the
# columns re-size to fit the new content.
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 16.333333333333336

# --- zh-cn detail sheet: Status (col C) + Latest Handoff Datetime (col H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-29 10:59:42"
$zhcn.Range("C1").EntireColumn.ColumnWidth = 16.333333333333336

# --- de-de detail sheet: Status (col C) + Latest Handoff Datetime (col H) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-29 10:59:46"
$dede.Range("C1").EntireColumn.ColumnWidth = 16.333333333333336
